$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.22658
$ws.Range("B3").Value = 0.22558
$ws.Range("B4").Value = 0.23725
$ws.Range("B5").Value = 0.28676
$ws.Range("B6").Value = 0.2802
$ws.Range("B7").Value = 0.2598
$ws.Range("B8").Value = 0.26462
$ws.Range("B9").Value = 0.24207
